$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates as described by the commit diff (price/volume refresh
# plus two coin-row swaps: rows 16/17 and rows 44/45).
# D-column "Price" values are plain text (European-style dotted numbers),
# so force Text number format before assigning to avoid Excel auto-converting
# them to numeric/date values.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.530.34'
$ws.Range('E2').Value = '  -0.48%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.867.78'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.70'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4550'
$ws.Range('E7').Value = '  -1.79%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3837'
$ws.Range('E8').Value = '  -1.88%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07828'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9888'
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.51'
$ws.Range('E11').Value = '  -3.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.848.87'
$ws.Range('E12').Value = '  +1.33%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.917'
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.643'
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06926'
$ws.Range('E15').Value = '  -0.15%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '86.46'
$ws.Range('E16').Value = '  -2.46%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.005'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009936'
$ws.Range('E18').Value = '  -1.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.68'
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.007'
$ws.Range('E20').Value = '  +0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.524.13'
$ws.Range('E21').Value = '  -0.39%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.254'
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.90'
$ws.Range('E23').Value = '  -1.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.083'
$ws.Range('E24').Value = '  -1.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.098.48'
$ws.Range('E25').Value = '  +2.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.17'
$ws.Range('E26').Value = '  -1.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.09'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.685'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '117.45'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.890'
$ws.Range('E30').Value = '  -5.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09277'
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9066'
$ws.Range('E32').Value = '  -3.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.277'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.320'
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.262'
$ws.Range('E35').Value = '  -2.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05654'
$ws.Range('E36').Value = '  -2.87%  '
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02037'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.620'
$ws.Range('E39').Value = '  -3.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5570'
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1765'
$ws.Range('E41').Value = '  -0.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.672'
$ws.Range('E42').Value = '  -2.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07143'
$ws.Range('E43').Value = '  -1.68%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5245'
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.48'
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.137'
$ws.Range('E46').Value = '  -4.32%  '
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.807'
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '111.54'
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.430'
$ws.Range('E50').Value = '  +3.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  -0.15%  '
